# Applies a data repull to column F (dSF) for a set of rows in Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value for column F
$updates = @{
    2  = -7
    3  = 0
    4  = -1
    8  = 0
    9  = -5
    13 = 1
    21 = 3
    22 = -1
    23 = 3
    25 = 4
    31 = 0
    36 = -6
    40 = -4
    44 = -2
    45 = 2
    48 = 5
    51 = 7
    55 = -2
    56 = -7
    57 = 0
    58 = -5
    59 = -3
    64 = 2
    66 = -8
    67 = -5
    69 = -9
    73 = -8
    77 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
